$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 'BTC'
$ws.Cells.Item(2, 3).Value = 'Bitcoin'
$ws.Cells.Item(2, 4).Value = 69360
$ws.Cells.Item(2, 5).Value = 1367250752088
$ws.Cells.Item(2, 6).Value = 5391088032
$ws.Cells.Item(2, 7).Value = 0.04864
$ws.Cells.Item(3, 2).Value = 'ETH'
$ws.Cells.Item(3, 3).Value = 'Ethereum'
$ws.Cells.Item(3, 4).Value = 3691.22
$ws.Cells.Item(3, 5).Value = 443837405773
$ws.Cells.Item(3, 6).Value = 6365866605
$ws.Cells.Item(3, 7).Value = 0.22743
$ws.Cells.Item(4, 2).Value = 'USDT'
$ws.Cells.Item(4, 3).Value = 'Tether'
$ws.Cells.Item(4, 4).Value = 0.999941
$ws.Cells.Item(4, 5).Value = 112446454286
$ws.Cells.Item(4, 6).Value = 20888854477
$ws.Cells.Item(4, 7).Value = 0.03079
$ws.Cells.Item(5, 2).Value = 'BNB'
$ws.Cells.Item(5, 3).Value = 'BNB'
$ws.Cells.Item(5, 4).Value = 680.4
$ws.Cells.Item(5, 5).Value = 104586807016
$ws.Cells.Item(5, 6).Value = 859051521
$ws.Cells.Item(5, 7).Value = -0.47033
$ws.Cells.Item(6, 2).Value = 'SOL'
$ws.Cells.Item(6, 3).Value = 'Solana'
$ws.Cells.Item(6, 4).Value = 159.44
$ws.Cells.Item(6, 5).Value = 73485795630
$ws.Cells.Item(6, 6).Value = 1462274716
$ws.Cells.Item(6, 7).Value = -0.47263
$ws.Cells.Item(7, 2).Value = 'STETH'
$ws.Cells.Item(7, 3).Value = 'Lido Staked Ether'
$ws.Cells.Item(7, 4).Value = 3690.87
$ws.Cells.Item(7, 5).Value = 35151542265
$ws.Cells.Item(7, 6).Value = 38833881
$ws.Cells.Item(7, 7).Value = 0.22836
$ws.Cells.Item(8, 2).Value = 'USDC'
$ws.Cells.Item(8, 3).Value = 'USDC'
$ws.Cells.Item(8, 4).Value = 1
$ws.Cells.Item(8, 5).Value = 32221765797
$ws.Cells.Item(8, 6).Value = 1699440036
$ws.Cells.Item(8, 7).Value = -0.00027
$ws.Cells.Item(9, 2).Value = 'XRP'
$ws.Cells.Item(9, 3).Value = 'XRP'
$ws.Cells.Item(9, 4).Value = 0.49439
$ws.Cells.Item(9, 5).Value = 27441349530
$ws.Cells.Item(9, 6).Value = 525066688
$ws.Cells.Item(9, 7).Value = 0.09417
$ws.Cells.Item(10, 2).Value = 'DOGE'
$ws.Cells.Item(10, 3).Value = 'Dogecoin'
$ws.Cells.Item(10, 4).Value = 0.146092
$ws.Cells.Item(10, 5).Value = 21150589907
$ws.Cells.Item(10, 6).Value = 708537293
$ws.Cells.Item(10, 7).Value = 0.28999
$ws.Cells.Item(11, 2).Value = 'TON'
$ws.Cells.Item(11, 3).Value = 'Toncoin'
$ws.Cells.Item(11, 4).Value = 7.13
$ws.Cells.Item(11, 5).Value = 17346977268
$ws.Cells.Item(11, 6).Value = 289865772
$ws.Cells.Item(11, 7).Value = -1.8493
$ws.Cells.Item(12, 2).Value = 'ADA'
$ws.Cells.Item(12, 3).Value = 'Cardano'
$ws.Cells.Item(12, 4).Value = 0.440274
$ws.Cells.Item(12, 5).Value = 15593014921
$ws.Cells.Item(12, 6).Value = 336163380
$ws.Cells.Item(12, 7).Value = 1.28801
$ws.Cells.Item(13, 2).Value = 'SHIB'
$ws.Cells.Item(13, 3).Value = 'Shiba Inu'
$ws.Cells.Item(13, 4).Value = 0.00002313
$ws.Cells.Item(13, 5).Value = 13639842151
$ws.Cells.Item(13, 6).Value = 423489431
$ws.Cells.Item(13, 7).Value = -1.26681
$ws.Cells.Item(14, 2).Value = 'AVAX'
$ws.Cells.Item(14, 3).Value = 'Avalanche'
$ws.Cells.Item(14, 4).Value = 32.29
$ws.Cells.Item(14, 5).Value = 12696144106
$ws.Cells.Item(14, 6).Value = 257207003
$ws.Cells.Item(14, 7).Value = -1.06275
$ws.Cells.Item(15, 2).Value = 'WBTC'
$ws.Cells.Item(15, 3).Value = 'Wrapped Bitcoin'
$ws.Cells.Item(15, 4).Value = 69402
$ws.Cells.Item(15, 5).Value = 10622473308
$ws.Cells.Item(15, 6).Value = 59821082
$ws.Cells.Item(15, 7).Value = 0.13685
$ws.Cells.Item(16, 2).Value = 'TRX'
$ws.Cells.Item(16, 3).Value = 'TRON'
$ws.Cells.Item(16, 4).Value = 0.116338
$ws.Cells.Item(16, 5).Value = 10161971502
$ws.Cells.Item(16, 6).Value = 372333833
$ws.Cells.Item(16, 7).Value = 3.03667
$ws.Cells.Item(17, 2).Value = 'LINK'
$ws.Cells.Item(17, 3).Value = 'Chainlink'
$ws.Cells.Item(17, 4).Value = 15.97
$ws.Cells.Item(17, 5).Value = 9378516167
$ws.Cells.Item(17, 6).Value = 270212396
$ws.Cells.Item(17, 7).Value = 0.72202
$ws.Cells.Item(18, 2).Value = 'BCH'
$ws.Cells.Item(18, 3).Value = 'Bitcoin Cash'
$ws.Cells.Item(18, 4).Value = 467.83
$ws.Cells.Item(18, 5).Value = 9229095789
$ws.Cells.Item(18, 6).Value = 195713396
$ws.Cells.Item(18, 7).Value = -0.6256699999999999
$ws.Cells.Item(19, 2).Value = 'DOT'
$ws.Cells.Item(19, 3).Value = 'Polkadot'
$ws.Cells.Item(19, 4).Value = 6.46
$ws.Cells.Item(19, 5).Value = 8876322274
$ws.Cells.Item(19, 6).Value = 193568163
$ws.Cells.Item(19, 7).Value = 0.1943
$ws.Cells.Item(20, 2).Value = 'UNI'
$ws.Cells.Item(20, 3).Value = 'Uniswap'
$ws.Cells.Item(20, 4).Value = 9.91
$ws.Cells.Item(20, 5).Value = 7470370046
$ws.Cells.Item(20, 6).Value = 210102584
$ws.Cells.Item(20, 7).Value = -0.05113
$ws.Cells.Item(21, 2).Value = 'NEAR'
$ws.Cells.Item(21, 3).Value = 'NEAR Protocol'
$ws.Cells.Item(21, 4).Value = 6.54
$ws.Cells.Item(21, 5).Value = 7091323441
$ws.Cells.Item(21, 6).Value = 239600691
$ws.Cells.Item(21, 7).Value = -2.7461
$ws.Cells.Item(22, 2).Value = 'MATIC'
$ws.Cells.Item(22, 3).Value = 'Polygon'
$ws.Cells.Item(22, 4).Value = 0.651784
$ws.Cells.Item(22, 5).Value = 6055759392
$ws.Cells.Item(22, 6).Value = 284662232
$ws.Cells.Item(22, 7).Value = -0.03113
$ws.Cells.Item(23, 2).Value = 'LTC'
$ws.Cells.Item(23, 3).Value = 'Litecoin'
$ws.Cells.Item(23, 4).Value = 80.13
$ws.Cells.Item(23, 5).Value = 5983680003
$ws.Cells.Item(23, 6).Value = 258747745
$ws.Cells.Item(23, 7).Value = 0.68152
$ws.Cells.Item(24, 2).Value = 'WEETH'
$ws.Cells.Item(24, 3).Value = 'Wrapped eETH'
$ws.Cells.Item(24, 4).Value = 3835.58
$ws.Cells.Item(24, 5).Value = 5538965554
$ws.Cells.Item(24, 6).Value = 29109773
$ws.Cells.Item(24, 7).Value = 0.19887
$ws.Cells.Item(25, 2).Value = 'LEO'
$ws.Cells.Item(25, 3).Value = 'LEO Token'
$ws.Cells.Item(25, 4).Value = 5.82
$ws.Cells.Item(25, 5).Value = 5394414915
$ws.Cells.Item(25, 6).Value = 1132797
$ws.Cells.Item(25, 7).Value = -2.43816
$ws.Cells.Item(26, 2).Value = 'DAI'
$ws.Cells.Item(26, 3).Value = 'Dai'
$ws.Cells.Item(26, 4).Value = 0.999629
$ws.Cells.Item(26, 5).Value = 5264485365
$ws.Cells.Item(26, 6).Value = 356103332
$ws.Cells.Item(26, 7).Value = 0.07213
$ws.Cells.Item(27, 2).Value = 'PEPE'
$ws.Cells.Item(27, 3).Value = 'Pepe'
$ws.Cells.Item(27, 4).Value = 0.00001229
$ws.Cells.Item(27, 5).Value = 5183671460
$ws.Cells.Item(27, 6).Value = 915086505
$ws.Cells.Item(27, 7).Value = -3.32052
$ws.Cells.Item(28, 2).Value = 'ICP'
$ws.Cells.Item(28, 3).Value = 'Internet Computer'
$ws.Cells.Item(28, 4).Value = 10.91
$ws.Cells.Item(28, 5).Value = 5073896453
$ws.Cells.Item(28, 6).Value = 112212226
$ws.Cells.Item(28, 7).Value = -1.46039
$ws.Cells.Item(29, 2).Value = 'FET'
$ws.Cells.Item(29, 3).Value = 'Fetch.ai'
$ws.Cells.Item(29, 4).Value = 1.73
$ws.Cells.Item(29, 5).Value = 4360497542
$ws.Cells.Item(29, 6).Value = 169817660
$ws.Cells.Item(29, 7).Value = -1.61457
$ws.Cells.Item(30, 2).Value = 'ETC'
$ws.Cells.Item(30, 3).Value = 'Ethereum Classic'
$ws.Cells.Item(30, 4).Value = 26.96
$ws.Cells.Item(30, 5).Value = 3975235754
$ws.Cells.Item(30, 6).Value = 119778644
$ws.Cells.Item(30, 7).Value = 0.86816
$ws.Cells.Item(31, 2).Value = 'KAS'
$ws.Cells.Item(31, 3).Value = 'Kaspa'
$ws.Cells.Item(31, 4).Value = 0.15797
$ws.Cells.Item(31, 5).Value = 3771816953
$ws.Cells.Item(31, 6).Value = 53432591
$ws.Cells.Item(31, 7).Value = -2.47027
$ws.Cells.Item(32, 2).Value = 'EZETH'
$ws.Cells.Item(32, 3).Value = 'Renzo Restaked ETH'
$ws.Cells.Item(32, 4).Value = 3680.42
$ws.Cells.Item(32, 5).Value = 3671852348
$ws.Cells.Item(32, 6).Value = 46702882
$ws.Cells.Item(32, 7).Value = 0.91612
$ws.Cells.Item(33, 2).Value = 'APT'
$ws.Cells.Item(33, 3).Value = 'Aptos'
$ws.Cells.Item(33, 4).Value = 8.31
$ws.Cells.Item(33, 5).Value = 3639911893
$ws.Cells.Item(33, 6).Value = 474979193
$ws.Cells.Item(33, 7).Value = 2.341
$ws.Cells.Item(34, 2).Value = 'RNDR'
$ws.Cells.Item(34, 3).Value = 'Render'
$ws.Cells.Item(34, 4).Value = 9.130000000000001
$ws.Cells.Item(34, 5).Value = 3552702522
$ws.Cells.Item(34, 6).Value = 117519541
$ws.Cells.Item(34, 7).Value = -0.00223
$ws.Cells.Item(35, 2).Value = 'FIL'
$ws.Cells.Item(35, 3).Value = 'Filecoin'
$ws.Cells.Item(35, 4).Value = 6.25
$ws.Cells.Item(35, 5).Value = 3528024040
$ws.Cells.Item(35, 6).Value = 279697274
$ws.Cells.Item(35, 7).Value = 2.25528
$ws.Cells.Item(36, 2).Value = 'USDE'
$ws.Cells.Item(36, 3).Value = 'Ethena USDe'
$ws.Cells.Item(36, 4).Value = 1.001
$ws.Cells.Item(36, 5).Value = 3358999868
$ws.Cells.Item(36, 6).Value = 86596221
$ws.Cells.Item(36, 7).Value = 0.08747000000000001
$ws.Cells.Item(37, 2).Value = 'FDUSD'
$ws.Cells.Item(37, 3).Value = 'First Digital USD'
$ws.Cells.Item(37, 4).Value = 0.999739
$ws.Cells.Item(37, 5).Value = 3277803320
$ws.Cells.Item(37, 6).Value = 2685370293
$ws.Cells.Item(37, 7).Value = -0.04643
$ws.Cells.Item(38, 2).Value = 'STX'
$ws.Cells.Item(38, 3).Value = 'Stacks'
$ws.Cells.Item(38, 4).Value = 2.23
$ws.Cells.Item(38, 5).Value = 3269110388
$ws.Cells.Item(38, 6).Value = 62105118
$ws.Cells.Item(38, 7).Value = -1.37302
$ws.Cells.Item(39, 2).Value = 'HBAR'
$ws.Cells.Item(39, 3).Value = 'Hedera'
$ws.Cells.Item(39, 4).Value = 0.09024799999999999
$ws.Cells.Item(39, 5).Value = 3228630561
$ws.Cells.Item(39, 6).Value = 54848058
$ws.Cells.Item(39, 7).Value = -0.59845
$ws.Cells.Item(40, 2).Value = 'XMR'
$ws.Cells.Item(40, 3).Value = 'Monero'
$ws.Cells.Item(40, 4).Value = 169.48
$ws.Cells.Item(40, 5).Value = 3126968415
$ws.Cells.Item(40, 6).Value = 92580678
$ws.Cells.Item(40, 7).Value = 2.55853
$ws.Cells.Item(41, 2).Value = 'MNT'
$ws.Cells.Item(41, 3).Value = 'Mantle'
$ws.Cells.Item(41, 4).Value = 0.941105
$ws.Cells.Item(41, 5).Value = 3075451282
$ws.Cells.Item(41, 6).Value = 103410213
$ws.Cells.Item(41, 7).Value = -0.49347
$ws.Cells.Item(42, 2).Value = 'ATOM'
$ws.Cells.Item(42, 3).Value = 'Cosmos Hub'
$ws.Cells.Item(42, 4).Value = 7.83
$ws.Cells.Item(42, 5).Value = 3058099471
$ws.Cells.Item(42, 6).Value = 151033583
$ws.Cells.Item(42, 7).Value = -0.8179999999999999
$ws.Cells.Item(43, 2).Value = 'CRO'
$ws.Cells.Item(43, 3).Value = 'Cronos'
$ws.Cells.Item(43, 4).Value = 0.111544
$ws.Cells.Item(43, 5).Value = 2992440260
$ws.Cells.Item(43, 6).Value = 8261352
$ws.Cells.Item(43, 7).Value = -0.90934
$ws.Cells.Item(44, 2).Value = 'IMX'
$ws.Cells.Item(44, 3).Value = 'Immutable'
$ws.Cells.Item(44, 4).Value = 1.98
$ws.Cells.Item(44, 5).Value = 2945446764
$ws.Cells.Item(44, 6).Value = 85315517
$ws.Cells.Item(44, 7).Value = -2.15733
$ws.Cells.Item(45, 2).Value = 'XLM'
$ws.Cells.Item(45, 3).Value = 'Stellar'
$ws.Cells.Item(45, 4).Value = 0.099231
$ws.Cells.Item(45, 5).Value = 2886348960
$ws.Cells.Item(45, 6).Value = 38091569
$ws.Cells.Item(45, 7).Value = 1.64986
$ws.Cells.Item(46, 2).Value = 'OKB'
$ws.Cells.Item(46, 3).Value = 'OKB'
$ws.Cells.Item(46, 4).Value = 47.11
$ws.Cells.Item(46, 5).Value = 2828276549
$ws.Cells.Item(46, 6).Value = 4516645
$ws.Cells.Item(46, 7).Value = -1.51584
$ws.Cells.Item(47, 2).Value = 'ARB'
$ws.Cells.Item(47, 3).Value = 'Arbitrum'
$ws.Cells.Item(47, 4).Value = 0.974952
$ws.Cells.Item(47, 5).Value = 2825777370
$ws.Cells.Item(47, 6).Value = 239652288
$ws.Cells.Item(47, 7).Value = 0.64466
$ws.Cells.Item(48, 2).Value = 'FLOKI'
$ws.Cells.Item(48, 3).Value = 'FLOKI'
$ws.Cells.Item(48, 4).Value = 0.0002792
$ws.Cells.Item(48, 5).Value = 2708506639
$ws.Cells.Item(48, 6).Value = 384552280
$ws.Cells.Item(48, 7).Value = 1.65058
$ws.Cells.Item(49, 2).Value = 'WIF'
$ws.Cells.Item(49, 3).Value = 'dogwifhat'
$ws.Cells.Item(49, 4).Value = 2.7
$ws.Cells.Item(49, 5).Value = 2697267573
$ws.Cells.Item(49, 6).Value = 360627832
$ws.Cells.Item(49, 7).Value = -0.99758
$ws.Cells.Item(50, 2).Value = 'INJ'
$ws.Cells.Item(50, 3).Value = 'Injective'
$ws.Cells.Item(50, 4).Value = 27.92
$ws.Cells.Item(50, 5).Value = 2684293840
$ws.Cells.Item(50, 6).Value = 199037533
$ws.Cells.Item(50, 7).Value = -4.71502
$ws.Cells.Item(51, 2).Value = 'SUI'
$ws.Cells.Item(51, 3).Value = 'Sui'
$ws.Cells.Item(51, 4).Value = 1.099
$ws.Cells.Item(51, 5).Value = 2669283318
$ws.Cells.Item(51, 6).Value = 307442437
$ws.Cells.Item(51, 7).Value = 0.40326
